$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 202.5
$ws.Cells.Item(4, 9).Value = 202.5
$ws.Cells.Item(4, 11).Value = 202.5
$ws.Cells.Item(4, 13).Value = -88.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 10807.167
$ws.Cells.Item(40, 9).Value = 12729.667
$ws.Cells.Item(40, 10).Value = 8884.666999999999
$ws.Cells.Item(40, 11).Value = 12729.667
$ws.Cells.Item(40, 12).Value = 8884.666999999999
$ws.Cells.Item(40, 13).Value = -12554.667
$ws.Cells.Item(40, 14).Value = -9234.666999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 7587.2
$ws.Cells.Item(74, 9).Value = 7587.2
$ws.Cells.Item(74, 11).Value = 7587.2
$ws.Cells.Item(74, 13).Value = -6651.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 7587.2
$ws.Cells.Item(77, 9).Value = 7587.2
$ws.Cells.Item(77, 11).Value = 37936
$ws.Cells.Item(77, 13).Value = -33256

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 7822.8335
$ws.Cells.Item(116, 9).Value = 32999
$ws.Cells.Item(116, 10).Value = 6341.8823
$ws.Cells.Item(116, 11).Value = 32999
$ws.Cells.Item(116, 12).Value = 6341.8823
$ws.Cells.Item(116, 13).Value = -29557
$ws.Cells.Item(116, 14).Value = -13225.8823

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 21451.824
$ws.Cells.Item(132, 9).Value = 3466.342
$ws.Cells.Item(132, 10).Value = 57422.79
$ws.Cells.Item(132, 11).Value = 10399.026
$ws.Cells.Item(132, 12).Value = 172268.37
$ws.Cells.Item(132, 13).Value = -7869.026
$ws.Cells.Item(132, 14).Value = -177328.37

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 1614.125
$ws.Cells.Item(135, 9).Value = 1164.7142
$ws.Cells.Item(135, 10).Value = 2243.3
$ws.Cells.Item(135, 11).Value = 10482.4278
$ws.Cells.Item(135, 12).Value = 20189.7
$ws.Cells.Item(135, 13).Value = -7947.427799999999
$ws.Cells.Item(135, 14).Value = -25259.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 5543
$ws.Cells.Item(137, 9).Value = 5218.3335
$ws.Cells.Item(137, 11).Value = 15655.0005
$ws.Cells.Item(137, 13).Value = -13105.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 955.9
$ws.Cells.Item(61, 9).Value = 955.2857
$ws.Cells.Item(61, 10).Value = 957.3333
$ws.Cells.Item(61, 11).Value = 955.2857
$ws.Cells.Item(61, 12).Value = 957.3333
$ws.Cells.Item(61, 13).Value = -743.2857
$ws.Cells.Item(61, 14).Value = -1381.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 1998
$ws.Cells.Item(63, 9).Value = 1816.091
$ws.Cells.Item(63, 10).Value = 2998.5
$ws.Cells.Item(63, 11).Value = 1816.091
$ws.Cells.Item(63, 12).Value = 2998.5
$ws.Cells.Item(63, 13).Value = -1130.091
$ws.Cells.Item(63, 14).Value = -4370.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 1998
$ws.Cells.Item(66, 9).Value = 1816.091
$ws.Cells.Item(66, 10).Value = 2998.5
$ws.Cells.Item(66, 11).Value = 9080.455
$ws.Cells.Item(66, 12).Value = 14992.5
$ws.Cells.Item(66, 13).Value = -5648.455
$ws.Cells.Item(66, 14).Value = -21856.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 5909.2
$ws.Cells.Item(102, 9).Value = 5685
$ws.Cells.Item(102, 10).Value = 6165.4287
$ws.Cells.Item(102, 11).Value = 5685
$ws.Cells.Item(102, 12).Value = 6165.4287
$ws.Cells.Item(102, 13).Value = -4063
$ws.Cells.Item(102, 14).Value = -9409.4287

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 8341.333000000001
$ws.Cells.Item(132, 9).Value = 4555.1816
$ws.Cells.Item(132, 11).Value = 13665.5448
$ws.Cells.Item(132, 13).Value = -11135.5448

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 955.9
$ws.Cells.Item(136, 9).Value = 955.2857
$ws.Cells.Item(136, 10).Value = 957.3333
$ws.Cells.Item(136, 11).Value = 2865.8571
$ws.Cells.Item(136, 12).Value = 2871.9999
$ws.Cells.Item(136, 13).Value = -315.8571000000002
$ws.Cells.Item(136, 14).Value = -7971.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 12843.272
$ws.Cells.Item(82, 9).Value = 4586.5557
$ws.Cells.Item(82, 10).Value = 49998.5
$ws.Cells.Item(82, 11).Value = 4586.5557
$ws.Cells.Item(82, 12).Value = 49998.5
$ws.Cells.Item(82, 13).Value = -4203.5557
$ws.Cells.Item(82, 14).Value = -50764.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(85, 8).Value = 12843.272
$ws.Cells.Item(85, 9).Value = 4586.5557
$ws.Cells.Item(85, 10).Value = 49998.5
$ws.Cells.Item(85, 11).Value = 4586.5557
$ws.Cells.Item(85, 12).Value = 49998.5
$ws.Cells.Item(85, 13).Value = -3260.5557
$ws.Cells.Item(85, 14).Value = -52650.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2540.3572
$ws.Cells.Item(134, 9).Value = 2540.3572
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 7621.071599999999
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 6766.3965
$ws.Cells.Item(58, 9).Value = 6797.775
$ws.Cells.Item(58, 10).Value = 6696.6665
$ws.Cells.Item(58, 11).Value = 6797.775
$ws.Cells.Item(58, 12).Value = 6696.6665
$ws.Cells.Item(58, 13).Value = -6594.775
$ws.Cells.Item(58, 14).Value = -7102.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 6766.3965
$ws.Cells.Item(136, 9).Value = 6797.775
$ws.Cells.Item(136, 10).Value = 6696.6665
$ws.Cells.Item(136, 11).Value = 20393.325
$ws.Cells.Item(136, 12).Value = 20089.9995
$ws.Cells.Item(136, 13).Value = -17843.325
$ws.Cells.Item(136, 14).Value = -25189.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(35, 8).Value = 1515.8889
$ws.Cells.Item(35, 9).Value = 400
$ws.Cells.Item(35, 10).Value = 2408.6
$ws.Cells.Item(35, 11).Value = 1200
$ws.Cells.Item(35, 12).Value = 7225.799999999999
$ws.Cells.Item(35, 13).Value = -912
$ws.Cells.Item(35, 14).Value = -7801.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8548.395500000001
$ws.Cells.Item(70, 9).Value = 8860.799999999999
$ws.Cells.Item(70, 11).Value = 8860.799999999999
$ws.Cells.Item(70, 13).Value = -8590.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 8548.395500000001
$ws.Cells.Item(73, 9).Value = 8860.799999999999
$ws.Cells.Item(73, 11).Value = 8860.799999999999
$ws.Cells.Item(73, 13).Value = -7924.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 6415.077
$ws.Cells.Item(132, 9).Value = 6280.2583
$ws.Cells.Item(132, 11).Value = 18840.7749
$ws.Cells.Item(132, 13).Value = -16310.7749

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(137, 8).Value = 100000
$ws.Cells.Item(137, 10).Value = 100000
$ws.Cells.Item(137, 12).Value = 100000
$ws.Cells.Item(137, 14).Value = -110200

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3972.3333
$ws.Cells.Item(122, 9).Value = 3842.1667
$ws.Cells.Item(122, 10).Value = 4232.6665
$ws.Cells.Item(122, 11).Value = 11526.5001
$ws.Cells.Item(122, 12).Value = 12697.9995
$ws.Cells.Item(122, 13).Value = -9076.500100000001
$ws.Cells.Item(122, 14).Value = -17597.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3657.2144
$ws.Cells.Item(132, 9).Value = 2886.1428
$ws.Cells.Item(132, 10).Value = 4428.2856
$ws.Cells.Item(132, 11).Value = 8658.428400000001
$ws.Cells.Item(132, 12).Value = 13284.8568
$ws.Cells.Item(132, 13).Value = -6128.428400000001
$ws.Cells.Item(132, 14).Value = -18344.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 11998.637
$ws.Cells.Item(136, 9).Value = 3198.6
$ws.Cells.Item(136, 10).Value = 99999
$ws.Cells.Item(136, 11).Value = 9595.799999999999
$ws.Cells.Item(136, 12).Value = 299997
$ws.Cells.Item(136, 13).Value = -7045.799999999999
$ws.Cells.Item(136, 14).Value = -305097

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 29330
$ws.Cells.Item(51, 9).Value = 26495
$ws.Cells.Item(51, 11).Value = 26495
$ws.Cells.Item(51, 13).Value = -25985

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1576.4242
$ws.Cells.Item(136, 9).Value = 1408.48
$ws.Cells.Item(136, 10).Value = 2101.25
$ws.Cells.Item(136, 11).Value = 4225.440000000001
$ws.Cells.Item(136, 12).Value = 6303.75
$ws.Cells.Item(136, 13).Value = -1675.440000000001
$ws.Cells.Item(136, 14).Value = -11403.75
